$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.165.46'
$ws.Range("E2").Value = '  -0.60%  '

$ws.Range("D3").Value = '1.588.36'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''211.52'
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("D6").Value = '''0.500'
$ws.Range("E6").Value = '  -0.79%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").Value = '''0.0600'
$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("D10").Value = '''19.01'
$ws.Range("E10").Value = '  -2.38%  '

$ws.Range("D11").Value = '''0.0843'
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").Value = '1.812.25'
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").Value = '1.581.37'
$ws.Range("E13").Value = '  -0.43%  '

$ws.Range("D14").Value = '''4.00'
$ws.Range("E14").Value = '  -1.81%  '

$ws.Range("D15").Value = '''0.510'
$ws.Range("E15").Value = '  -1.69%  '

$ws.Range("D16").Value = '''63.46'
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("D17").Value = '26.175.21'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '0.0₃0722'
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").Value = '''7.37'
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("D20").Value = '''213.49'
$ws.Range("E20").Value = '  +1.23%  '

$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("D23").Value = '''8.91'
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("E24").Value = '  -1.60%  '

$ws.Range("D25").Value = '''144.06'
$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").Value = '''6.93'
$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("E28").Value = '  -1.30%  '

$ws.Range("D29").Value = '''15.00'
$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("D30").Value = '''0.0492'
$ws.Range("E30").Value = '  -2.73%  '

$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '''3.14'
$ws.Range("E32").Value = '  -2.19%  '

$ws.Range("D33").Value = '1.410.12'
$ws.Range("E33").Value = '  +8.38%  '

$ws.Range("D34").Value = '''2.94'
$ws.Range("E34").Value = '  -1.73%  '

$ws.Range("E35").Value = '  -0.73%  '

$ws.Range("E36").Value = '  -1.48%  '

$ws.Range("D37").Value = '''0.584'
$ws.Range("E37").Value = '  -4.42%  '

$ws.Range("D38").Value = '''0.0165'
$ws.Range("E38").Value = '  -1.76%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '''0.819'
$ws.Range("E39").Value = '  +1.58%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''5.87'
$ws.Range("E40").Value = '  +4.63%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").Value = '''0.945'
$ws.Range("E42").Value = '  -13.47%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.13'
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.763'
$ws.Range("E44").Value = '  -0.65%  '

$ws.Range("D45").Value = '1.723.06'
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").Value = '''60.82'
$ws.Range("E46").Value = '  -2.58%  '

$ws.Range("D47").Value = '''85.52'
$ws.Range("E47").Value = '  -2.56%  '

$ws.Range("E48").Value = '  -1.08%  '

$ws.Range("D49").Value = '''1.48'
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").Value = '''0.0501'
$ws.Range("E50").Value = '  -0.83%  '

$ws.Range("D51").Value = '''0.0953'
$ws.Range("E51").Value = '  -2.81%  '
